$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rows = @(
  @{row=2; A="ECs"; B="Inhba"; C="Acvr2b"; D="ECs"; E=2; F=0.6666666666666666; G=3.675031333333333; H=11.025094; I=0.2032371147293133; J=0.2032371147293133; K=3; L=1; M=0.5292956666666667; N=1.587887; O=0.135651968140022; P=0.1356519681400219; Q=1.945178159597556; R=17.506603436378; S=0.02756951461213079; T=0.02756951461213079}
  @{row=3; A="ECs"; B="Inhba"; C="Acvr2b"; D="FAPs"; E=2; F=0.6666666666666666; G=3.675031333333333; H=11.025094; I=0.2032371147293133; J=0.2032371147293133; K=3; L=1; M=1.362890666666667; N=4.088672000000001; O=0.3492921120199358; P=0.3492921120199358; Q=5.008665903907557; R=45.07799313516801; S=0.07098912104463985; T=0.07098912104463985}
  @{row=4; A="ECs"; B="Inhba"; C="Acvr2b"; D="sCs"; E=2; F=0.6666666666666666; G=3.675031333333333; H=11.025094; I=0.2032371147293133; J=0.2032371147293133; K=3; L=1; M=2.009678666666666; N=6.029036; O=0.5150559198400423; P=0.5150559198400423; Q=7.385632069931554; R=66.47068862938399; S=0.1046784790725427; T=0.1046784790725427}
  @{row=5; A="FAPs"; B="Inhba"; C="Acvr2b"; D="ECs"; E=3; F=1; G=10.108494; H=30.325482; I=0.5590213983169419; J=0.5590213983169419; K=3; L=1; M=0.5292956666666667; N=1.587887; O=0.135651968140022; P=0.1356519681400219; Q=5.350382070726001; R=48.15343863653401; S=0.07583235291408032; T=0.0758323529140803}
  @{row=6; A="FAPs"; B="Inhba"; C="Acvr2b"; D="FAPs"; E=3; F=1; G=10.108494; H=30.325482; I=0.5590213983169419; J=0.5590213983169419; K=3; L=1; M=1.362890666666667; N=4.088672000000001; O=0.3492921120199358; P=0.3492921120199358; Q=13.776772126656; R=123.990949139904; S=0.1952617648824624; T=0.1952617648824624}
  @{row=7; A="FAPs"; B="Inhba"; C="Acvr2b"; D="sCs"; E=3; F=1; G=10.108494; H=30.325482; I=0.5590213983169419; J=0.5590213983169419; K=3; L=1; M=2.009678666666666; N=6.029036; O=0.5150559198400423; P=0.5150559198400423; Q=20.314824743928; R=182.833422695352; S=0.2879272805203992; T=0.2879272805203992}
  @{row=8; A="sCs"; B="Inhba"; C="Acvr2b"; D="ECs"; E=3; F=1; G=4.298956; H=12.896868; I=0.2377414869537448; J=0.2377414869537448; K=3; L=1; M=0.5292956666666667; N=1.587887; O=0.135651968140022; P=0.1356519681400219; Q=2.275418781990667; R=20.478769037916; S=0.03225010061381083; T=0.03225010061381083}
  @{row=9; A="sCs"; B="Inhba"; C="Acvr2b"; D="FAPs"; E=3; F=1; G=4.298956; H=12.896868; I=0.2377414869537448; J=0.2377414869537448; K=3; L=1; M=1.362890666666667; N=4.088672000000001; O=0.3492921120199358; P=0.3492921120199358; Q=5.859007008810669; R=52.73106307929601; S=0.08304122609283353; T=0.08304122609283353}
  @{row=10; A="sCs"; B="Inhba"; C="Acvr2b"; D="sCs"; E=3; F=1; G=4.298956; H=12.896868; I=0.2377414869537448; J=0.2377414869537448; K=3; L=1; M=2.009678666666666; N=6.029036; O=0.5150559198400423; P=0.5150559198400423; Q=8.639520162138666; R=77.75568145924801; S=0.1224501602471004; T=0.1224501602471004}
)
foreach ($r in $rows) {
  $ws.Cells.Item($r.row, 1).Value = $r.A
  $ws.Cells.Item($r.row, 2).Value = $r.B
  $ws.Cells.Item($r.row, 3).Value = $r.C
  $ws.Cells.Item($r.row, 4).Value = $r.D
  $ws.Cells.Item($r.row, 5).Value = $r.E
  $ws.Cells.Item($r.row, 6).Value = $r.F
  $ws.Cells.Item($r.row, 7).Value = $r.G
  $ws.Cells.Item($r.row, 8).Value = $r.H
  $ws.Cells.Item($r.row, 9).Value = $r.I
  $ws.Cells.Item($r.row, 10).Value = $r.J
  $ws.Cells.Item($r.row, 11).Value = $r.K
  $ws.Cells.Item($r.row, 12).Value = $r.L
  $ws.Cells.Item($r.row, 13).Value = $r.M
  $ws.Cells.Item($r.row, 14).Value = $r.N
  $ws.Cells.Item($r.row, 15).Value = $r.O
  $ws.Cells.Item($r.row, 16).Value = $r.P
  $ws.Cells.Item($r.row, 17).Value = $r.Q
  $ws.Cells.Item($r.row, 18).Value = $r.R
  $ws.Cells.Item($r.row, 19).Value = $r.S
  $ws.Cells.Item($r.row, 20).Value = $r.T
}
